$wb = $excel.ActiveWorkbook

$wsIncome  = $wb.Worksheets.Item("NKE Income Annual")
$wsBalance = $wb.Worksheets.Item("NKE Balance Annual")
$wsCash    = $wb.Worksheets.Item("NKE Cash Annual")

# --- NKE Income Annual: drop the "Tax Rate for Calcs" block (rows 157-161) ---
$wsIncome.Rows("157:161").Delete()

# --- NKE Balance Annual: insert a new "TTM" header row above the first date row ---
$wsBalance.Rows("2:2").Insert()
$wsBalance.Range("A2").Value = "TTM"

# --- Update each sheet's saved selection / active cell ---
$wsIncome.Activate()
$wsIncome.Range("A143").Select()

$wsCash.Activate()
$wsCash.Range("D69").Select()

# Leave the workbook on the Balance sheet (matches the saved tabSelected/activeTab)
$wsBalance.Activate()
$wsBalance.Range("C7").Select()
